# caravel-Nucleo-v2-M.2 BOM: remove the "FlexyPins" row (row 35) left over
# from the previous iteration, leaving only the empty, still-formatted D35
# cell behind (matches the target OOXML exactly), and move the active
# selection to C38.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the hyperlink attached to D35 ("FlexyPins" -> lectronz.com link)
# before clearing the cell so the relationship is cleaned up too.
$ws.Range("D35").Hyperlinks.Delete()

# Fully clear A35:C35 and E35 (contents + formatting) so they disappear
# from the sheet entirely, matching the diff (no <c> emitted for them).
$ws.Range("A35").Clear()
$ws.Range("B35").Clear()
$ws.Range("C35").Clear()
$ws.Range("E35").Clear()

# D35 keeps its existing format/style (s="8") but loses its value/text.
$ws.Range("D35").ClearContents()

# Update the saved selection to match the author's final cursor position.
$ws.Range("C38").Select()
